$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 3635
$ws.Range("F5").Value = 3635
$ws.Range("F6").Value = 270
$ws.Range("F7").Value = 5162
$ws.Range("F8").Value = 547
$ws.Range("F9").Value = 377
$ws.Range("F10").Value = 204
$ws.Range("F11").Value = 703
$ws.Range("F13").Value = 102
$ws.Range("F14").Value = 37
$ws.Range("F15").Value = 712
$ws.Range("F16").Value = 323
$ws.Range("F17").Value = 38
$ws.Range("F19").Value = 160
$ws.Range("F21").Value = 363
$ws.Range("F22").Value = 4946
$ws.Range("F25").Value = 13
$ws.Range("F26").Value = 6071
$ws.Range("F29").Value = 3231
$ws.Range("F31").Value = 719
$ws.Range("F32").Value = 4448
$ws.Range("F34").Value = 126
$ws.Range("F35").Value = 143
$ws.Range("F36").Value = 1056
$ws.Range("F40").Value = 883
$ws.Range("F41").Value = 1038
$ws.Range("F42").Value = 2037

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 48
$ws.Range("F3").Value = 27
$ws.Range("F5").Value = 60

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 1127

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1127
$ws.Range("F7").Value = 3635
$ws.Range("F8").Value = 3635
$ws.Range("F9").Value = 270
$ws.Range("F10").Value = 5162
$ws.Range("F11").Value = 547
$ws.Range("F12").Value = 377
$ws.Range("F13").Value = 204
$ws.Range("F14").Value = 703
$ws.Range("F16").Value = 102
$ws.Range("F17").Value = 37
$ws.Range("F18").Value = 712
$ws.Range("F19").Value = 323
$ws.Range("F20").Value = 38
$ws.Range("F21").Value = 48
$ws.Range("F23").Value = 160
$ws.Range("F25").Value = 363
$ws.Range("F26").Value = 4946
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 6071
$ws.Range("F33").Value = 3231
$ws.Range("F35").Value = 719
$ws.Range("F36").Value = 4448
$ws.Range("F38").Value = 27
$ws.Range("F39").Value = 126
$ws.Range("F40").Value = 143
$ws.Range("F41").Value = 1056
$ws.Range("F45").Value = 883
$ws.Range("F46").Value = 1038
$ws.Range("F48").Value = 2037
$ws.Range("F50").Value = 60
